$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.919.04'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.828.43'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.58'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  +2.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.01'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  +4.80%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0990'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '2.092.91'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.845.27'
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.31'
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.668'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.64'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '34.868.59'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.60'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.92'
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.19'
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.66'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '173.91'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.76'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  +3.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.35'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  -6.32%  '
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0551'
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  +2.68%  '
$ws.Range("E35").Value = '  +5.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("E36").Value = '  +11.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.699'
$ws.Range("E37").Value = '  +3.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '92.32'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '1.338.71'
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.52'
$ws.Range("E42").Value = '  -2.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.25'
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.75'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0523'
$ws.Range("E47").Value = '  +1.88%  '
$ws.Range("D48").Value = '2.008.62'
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.86'
$ws.Range("E51").Value = '  -1.34%  '
